# Swap the presentation's applied theme from the custom "Integral" (Red
# Violet) design back to the stock Office "Office Theme" design - i.e. what
# happens when a new design is picked from the Design tab: the slide
# master's theme colours are replaced with the default Office palette.

function Convert-HexToRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$theme = $master.Theme

# Standard Office theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$colorScheme = $theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = Convert-HexToRGB $officeColors[$i - 1]
}
